$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name and title text to reflect new "through" date
$ws.Name = "Through 2022-08-20"
$ws.Range("A9").Value = "August (through 08-20)"

# Update August row (row 9) values
$ws.Range("B9").Value = 22
$ws.Range("C9").Value = 46
$ws.Range("D9").Value = 53
$ws.Range("G9").Value = 124
$ws.Range("H9").Value = 106
$ws.Range("I9").Value = 119

# Update Total row (row 10) values
$ws.Range("B10").Value = 184
$ws.Range("C10").Value = 348
$ws.Range("D10").Value = 518
$ws.Range("G10").Value = 745
$ws.Range("H10").Value = 1016
$ws.Range("I10").Value = 1090

$wb.Save()
